# "modificação na estrutura das pastas na pasta log_paste"
#
# Slide 13 ("Classe: Database") gets its connection-credentials callout
# reworked:
#   - the separator line is nudged up
#   - the "Name" bullet is renamed to "Password"
#   - the lower bullet list gains four new method names ahead of the
#     existing BuscaRegistro / InsereRegistro entries, and the textbox
#     is repositioned to stay visually centred as it grows

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

# 1) "Name" -> "Password" in the attribute bullet list (CaixaDeTexto 5)
$attrsShape = $s.Shapes.Item(5)
$attrsRange = $attrsShape.TextFrame.TextRange
$attrsRange.Text = $attrsRange.Text.Replace("Name", "Password")

# 2) Prepend the four new method bullets to CaixaDeTexto 6, ahead of the
#    existing BuscaRegistro / InsereRegistro bullets (insert in reverse
#    order so the final reading order is correct).
$methodsShape = $s.Shapes.Item(6)
$methodsRange = $methodsShape.TextFrame.TextRange
[void]$methodsRange.InsertBefore("DesconectarCursor`r")
[void]$methodsRange.InsertBefore("ConectarCursor`r")
[void]$methodsRange.InsertBefore("DesconectarBase`r")
[void]$methodsRange.InsertBefore("ConectarBase`r")

# The textbox auto-fits its height to the new content (spAutoFit); nudge
# its top up so it stays vertically centred against the rounded
# rectangle behind it, matching the author's manual re-layout.
$methodsShape.Top = 368.79252625984253

# 3) The connector line separating the two bullet lists moves up to sit
#    above the now-taller lower textbox.
$connectorShape = $s.Shapes.Item(3)
$connectorShape.Top = 346.9096062992126
